# WG Number excel table: add the missing WG Numbers (N9524 / N9525) for the
# new "extruded_structure_cross_section" ARM/MIM EXPRESS module rows, and
# drop the now-obsolete "modified"/"new" status markers in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WG NB")

# ---------------------------------------------------------------------
# 1) Column G ("status") no longer applies to rows 42-48 ("modified")
#    once the sheet is cleaned up - clear those cells entirely.
# ---------------------------------------------------------------------
$ws.Range("G42:G48").ClearContents()

# ---------------------------------------------------------------------
# 2) Row 49 - was a placeholder "N" / "new" row for the ARM EXPRESS
#    schema of the new module. Fill in the real WG number (N9524),
#    update the request date, restyle A/F to match the rest of the
#    table, and clear the now-unused status cell (keep it blank).
# ---------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A49").PasteSpecial(-4122)
$ws.Range("A49").Value = "N9524"

$ws.Range("F2").Copy()
$ws.Range("F49").PasteSpecial(-4122)
$ws.Range("F49").Value = 42817

$ws.Range("G49").ClearContents()

# ---------------------------------------------------------------------
# 3) Row 50 - same treatment for the MIM EXPRESS schema row (N9525).
# ---------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A50").PasteSpecial(-4122)
$ws.Range("A50").Value = "N9525"

$ws.Range("F2").Copy()
$ws.Range("F50").PasteSpecial(-4122)
$ws.Range("F50").Value = 42817

$ws.Range("G50").ClearContents()

# ---------------------------------------------------------------------
# 4) Restore the usual clipboard/marching-ants state and refresh the
#    window scroll position / selection to match where the author left
#    off editing (bottom of the now-longer table).
# ---------------------------------------------------------------------
$excel.CutCopyMode = $false | Out-Null

$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E53").Select() | Out-Null
